$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, copy the existing row 4 data down into a new row 5 (this preserves
# the old record: Femacal de La Calera, 10-kilo tray, Los Andes, 44874).
$ws.Range("A5").Value = $ws.Range("A4").Value()
$ws.Range("B5").Value = $ws.Range("B4").Value()
$ws.Range("C5").Value = $ws.Range("C4").Value()
$ws.Range("D5").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("D5").Value = $ws.Range("D4").Value()
$ws.Range("E5").Value = $ws.Range("E4").Value()
$ws.Range("F5").Value = $ws.Range("F4").Value()
$ws.Range("G5").Value = $ws.Range("G4").Value()
$ws.Range("H5").Value = $ws.Range("H4").Value()
$ws.Range("I5").Value = $ws.Range("I4").Value()
$ws.Range("J5").Value = $ws.Range("J4").Value()
$ws.Range("K5").Value = $ws.Range("K4").Value()
$ws.Range("L5").Value = $ws.Range("L4").Value()
$ws.Range("M5").Value = $ws.Range("M4").Value()
$ws.Range("N5").Value = $ws.Range("N4").Value()
$ws.Range("O5").Value = $ws.Range("O4").Value()
$ws.Range("P5").Value = $ws.Range("P4").Value()
$ws.Range("Q5").Value = $ws.Range("Q4").Value()
$ws.Range("R5").Value = $ws.Range("R4").Value()
$ws.Range("S5").Value = $ws.Range("S4").Value()
$ws.Range("T5").Value = $ws.Range("T4").Value()

# Now update row 4 with the new weekly values.
$ws.Range("D4").Value = 45222
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 1500
